$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match rows to append (Indice, pais, torneio, temporada, data_partida, home, home_ft_gols,
# away, away_ft_gols, home_opening_odds, home_opening_data_hora, home_closing_odds, home_closing_data_hora,
# draw_opening_odds, draw_opening_data_hora, draw_closing_odds, draw_closing_data_hora,
# away_opening_odds, away_opening_data_hora, away_closing_odds, away_closing_data_hora, url_partida)
$rows = @(
    @{
        Indice = 184
        data_partida = 45294.70833333334
        home = "Granada CF"
        home_ft_gols = 2
        away = "Cadiz CF"
        away_ft_gols = 0
        home_opening_odds = 1.78
        home_opening_data_hora = "17/12/2024 18:03"
        home_closing_odds = 2.25
        home_closing_data_hora = "03/01/2024 16:57"
        draw_opening_odds = 3.61
        draw_opening_data_hora = "17/12/2024 18:03"
        draw_closing_odds = 3.12
        draw_closing_data_hora = "03/01/2024 16:59"
        away_opening_odds = 4.41
        away_opening_data_hora = "17/12/2024 18:03"
        away_closing_odds = 3.82
        away_closing_data_hora = "03/01/2024 16:59"
        url_partida = "https://www.betexplorer.com/football/spain/laliga/granada-cf-cadiz/rVOyfxQG/"
    },
    @{
        Indice = 185
        data_partida = 45294.80208333334
        home = "Celta Vigo"
        home_ft_gols = 2
        away = "Betis"
        away_ft_gols = 1
        home_opening_odds = 1.95
        home_opening_data_hora = "17/12/2024 18:03"
        home_closing_odds = 2.21
        home_closing_data_hora = "03/01/2024 19:14"
        draw_opening_odds = 3.47
        draw_opening_data_hora = "17/12/2024 18:03"
        draw_closing_odds = 3.31
        draw_closing_data_hora = "03/01/2024 19:14"
        away_opening_odds = 3.77
        away_opening_data_hora = "17/12/2024 18:03"
        away_closing_odds = 3.67
        away_closing_data_hora = "03/01/2024 19:14"
        url_partida = "https://www.betexplorer.com/football/spain/laliga/celta-vigo-betis/URKXfIAA/"
    },
    @{
        Indice = 186
        data_partida = 45294.80208333334
        home = "Real Madrid"
        home_ft_gols = 1
        away = "Mallorca"
        away_ft_gols = 0
        home_opening_odds = 1.17
        home_opening_data_hora = "17/12/2024 18:03"
        home_closing_odds = 1.22
        home_closing_data_hora = "03/01/2024 19:10"
        draw_opening_odds = 6.63
        draw_opening_data_hora = "17/12/2024 18:03"
        draw_closing_odds = 6.75
        draw_closing_data_hora = "03/01/2024 19:14"
        away_opening_odds = 12.67
        away_opening_data_hora = "17/12/2024 18:03"
        away_closing_odds = 13.93
        away_closing_data_hora = "03/01/2024 19:14"
        url_partida = "https://www.betexplorer.com/football/spain/laliga/real-madrid-mallorca/xhOugduN/"
    }
)

$startRow = 185
$pais = "spain"
$torneio = "laliga"
$temporada = "2023-2024"

$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Indice

    $ws.Cells.Item($r, 2).Value = $pais
    $ws.Cells.Item($r, 3).Value = $torneio
    $ws.Cells.Item($r, 4).Value = $temporada

    $ws.Cells.Item($r, 5).Value = $row.data_partida

    $ws.Cells.Item($r, 6).Value = $row.home
    $ws.Cells.Item($r, 7).Value = $row.home_ft_gols
    $ws.Cells.Item($r, 8).Value = $row.away
    $ws.Cells.Item($r, 9).Value = $row.away_ft_gols

    $ws.Cells.Item($r, 10).Value = $row.home_opening_odds
    $ws.Cells.Item($r, 11).Value = $row.home_opening_data_hora
    $ws.Cells.Item($r, 12).Value = $row.home_closing_odds
    $ws.Cells.Item($r, 13).Value = $row.home_closing_data_hora

    $ws.Cells.Item($r, 14).Value = $row.draw_opening_odds
    $ws.Cells.Item($r, 15).Value = $row.draw_opening_data_hora
    $ws.Cells.Item($r, 16).Value = $row.draw_closing_odds
    $ws.Cells.Item($r, 17).Value = $row.draw_closing_data_hora

    $ws.Cells.Item($r, 18).Value = $row.away_opening_odds
    $ws.Cells.Item($r, 19).Value = $row.away_opening_data_hora
    $ws.Cells.Item($r, 20).Value = $row.away_closing_odds
    $ws.Cells.Item($r, 21).Value = $row.away_closing_data_hora

    $ws.Cells.Item($r, 22).Value = $row.url_partida

    $r = $r + 1
}

$lastRow = $startRow + $rows.Count - 1

# Match the formatting used throughout the table: column A (Indice) uses the
# bold/bordered/centered style, column E (data_partida) uses the
# YYYY-MM-DD HH:MM:SS number format. Copy formats from the row directly above
# the new block so the new rows pick up the existing shared styles exactly.
$ws.Cells.Item($startRow - 1, 1).Copy()
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($lastRow, 1)).PasteSpecial(-4122)

$ws.Cells.Item($startRow - 1, 5).Copy()
$ws.Range($ws.Cells.Item($startRow, 5), $ws.Cells.Item($lastRow, 5)).PasteSpecial(-4122)

$excel.CutCopyMode = 0
